# Apply the "Updated cryptos list" refresh (prices + 1h volume deltas,
# plus the VeChain/dogwifhat rank swap in rows 50-51) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "62.244.24"
$ws.Range("E2").Value = "  +1.96%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "3.430.21"
$ws.Range("E3").Value = "  +1.29%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  +0.03%  "

# Row 5: BNB
$ws.Range("D5").Value = "'579.93"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.33%  "

# Row 6: Solana
$ws.Range("D6").Value = "'145.69"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.02%  "

# Row 7: USDC
$ws.Range("E7").Value = "  +0.02%  "

# Row 8: XRP
$ws.Range("E8").Value = "  +0.19%  "

# Row 9: Toncoin
$ws.Range("D9").Value = "'7.61"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.70%  "

# Row 10: Dogecoin
$ws.Range("E10").Value = "  +1.12%  "

# Row 11: Cardano
$ws.Range("D11").Value = "'0.388"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.02%  "

# Row 12: WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "4.015.50"
$ws.Range("E12").Value = "  +1.21%  "

# Row 13: Avalanche
$ws.Range("D13").Value = "'28.90"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +3.60%  "

# Row 14: TRON
$ws.Range("E14").Value = "  -0.78%  "

# Row 15: WrappedEther
$ws.Range("D15").Value = "3.428.29"
$ws.Range("E15").Value = "  +1.08%  "

# Row 16: ShibaInu
$ws.Range("E16").Value = "  +0.31%  "

# Row 17: WrappedBTC
$ws.Range("D17").Value = "62.221.42"
$ws.Range("E17").Value = "  +1.78%  "

# Row 18: Polkadot
$ws.Range("D18").Value = "'6.21"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.46%  "

# Row 19: Chainlink
$ws.Range("D19").Value = "'14.09"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.94%  "

# Row 20: Uniswap
$ws.Range("D20").Value = "'9.23"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.01%  "

# Row 21: BitcoinCash
$ws.Range("D21").Value = "'394.18"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.45%  "

# Row 22: Litecoin
$ws.Range("D22").Value = "'75.00"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.50%  "

# Row 23: Polygon
$ws.Range("D23").Value = "'0.556"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.42%  "

# Row 24: Dai
$ws.Range("E24").Value = "  +0.03%  "

# Row 25: PEPE
$ws.Range("D25").Value = "'0.0000117"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.97%  "

# Row 26: WrappedeETH
$ws.Range("D26").Value = "3.566.98"
$ws.Range("E26").Value = "  +1.26%  "

# Row 27: Kaspa
$ws.Range("E27").Value = "  +0.60%  "

# Row 28: RenderToken
$ws.Range("D28").Value = "'7.55"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +4.37%  "

# Row 29: Binance-PegBSC-USD
$ws.Range("E29").Value = "  +0.01%  "

# Row 30: InternetComputer(DFINITY)
$ws.Range("D30").Value = "'8.05"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.06%  "

# Row 31: PancakeSwap
$ws.Range("E31").Value = "  +0.61%  "

# Row 32: Fetch.AI
$ws.Range("E32").Value = "  +3.05%  "

# Row 33: USDe
$ws.Range("E33").Value = "  -0.03%  "

# Row 34: EthereumClassic
$ws.Range("D34").Value = "'23.65"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.36%  "

# Row 35: NEARProtocol
$ws.Range("D35").Value = "'5.30"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +5.94%  "

# Row 36: Aptos
$ws.Range("D36").Value = "'7.02"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.81%  "

# Row 37: Monero
$ws.Range("D37").Value = "'167.77"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.26%  "

# Row 38: ImmutableX
$ws.Range("E38").Value = "  +4.18%  "

# Row 39: RenzoRestakedETH
$ws.Range("D39").Value = "3.462.30"
$ws.Range("E39").Value = "  +1.22%  "

# Row 40: EnergySwap
$ws.Range("D40").Value = "'28.75"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +8.00%  "

# Row 41: Hedera
$ws.Range("D41").Value = "'0.0755"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.67%  "

# Row 42: Mantle
$ws.Range("D42").Value = "'0.790"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.64%  "

# Row 43: Filecoin
$ws.Range("E43").Value = "  +1.91%  "

# Row 44: Stacks
$ws.Range("D44").Value = "'1.69"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.92%  "

# Row 45: ONDO
$ws.Range("D45").Value = "'1.18"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +4.77%  "

# Row 46: Maker
$ws.Range("D46").Value = "2.515.73"
$ws.Range("E46").Value = "  +2.26%  "

# Row 47: InjectiveProtocol
$ws.Range("D47").Value = "'23.16"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.15%  "

# Row 48: Cosmos
$ws.Range("D48").Value = "'6.68"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.36%  "

# Row 49: FirstDigitalUSD
$ws.Range("D49").Value = "'0.999"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.06%  "

# Row 50: dogwifhat
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0265"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.64%  "

# Row 51: VeChain
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").Value = "'2.14"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.63%  "
